# Automatische test-sync: 2025-08-04 20:14:50
# Adds the new test-mail #3 ("Kun jij dit afhandelen?") row to the Logs sheet,
# extends the conditional formatting ranges to cover it, and refreshes the
# Dashboard category counts / chart source rows to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Append the new log row (row 9) ---------------------------------------
$ws.Range("A9").Value = "Kun jij dit afhandelen?"
$ws.Range("B9").Value = "mailmind.test@zohomail.eu"
$ws.Range("C9").Value = "Testmail #3: Kun jij dit afhandelen?"
$ws.Range("D9").Value = "Planning / Afspraak"
$ws.Range("E9").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Range("F9").Value = "2025-08-04 20:14:02"
$ws.Range("G9").Value = "Ja"
$ws.Range("H9").Value = "Ja"
$ws.Range("I9").Value = "Nee"
$ws.Range("J9").Value = "Nee"

# --- Extend the conditional formatting ranges down to row 9 ---------------
$ws.Range("D2:D8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D9"))
$ws.Range("G2:G8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G9"))
$ws.Range("H2:H8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H9"))
$ws.Range("I2:I8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I9"))
$ws.Range("J2:J8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J9"))

# --- Refresh the Dashboard summary table -----------------------------------
# New counts: Opvolging/Status=3 (unchanged), Planning/Afspraak=3 (+1 new
# mail), Retour/Terugbetaling=2 (unchanged) - rows 3 & 4 swap order/labels
# so the table stays sorted by descending count.
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Range("A3").Value = "Planning / Afspraak"
$ws2.Range("B3").Value = 3
$ws2.Range("A4").Value = "Retour / Terugbetaling"
$ws2.Range("B4").Value = 2
